# Add two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# formatting of the existing header/data columns, and populate the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy style from H1 (bold, centered, bordered) onto I1/J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

# Data values for rows 2..18 -> columns I (col 9) and J (col 10)
$data = @{
    2  = @(2, 5)
    3  = @(1, 3)
    4  = @(3, 6)
    5  = @(1, 6)
    6  = @(2, 6)
    7  = @(1, 5)
    8  = @(1, 6)
    9  = @(1, 6)
    10 = @(1, 6)
    11 = @(2, 8)
    12 = @(1, 5)
    13 = @(2, 6)
    14 = @(2, 6)
    15 = @(5, 5)
    16 = @(8, 9)
    17 = @(7, 9)
    18 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
